# Delete the entire row that contains the "Calle Camilo Cienfuegos..." address
# (shared string index 32), which shifts all subsequent rows up by one and
# removes the now-unused shared string / redundant cell style automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$ws.Rows.Item(38).Delete()
$ws.Range("A28").Style = "Normal 2"
